$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the edited cells to remain plain text (matching the source's inlineStr
# cells) so Excel does not auto-convert numeric-looking strings (e.g. "19.50")
# into numbers and strip significant trailing zeros / thousand-dot separators.
$editRange = $ws.Range("D2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = '72.404.71'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '2.660.33'
$ws.Range("E3").Value = '  +0.76%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '596.93'
$ws.Range("E5").Value = '  -1.45%  '
$ws.Range("D6").Value = '175.39'
$ws.Range("E6").Value = '  -2.45%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("D9").Value = '2.658.00'
$ws.Range("E9").Value = '  +0.69%  '
$ws.Range("E10").Value = '  -3.34%  '
$ws.Range("E11").Value = '  +1.99%  '
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("E13").Value = '  -0.92%  '
$ws.Range("D14").Value = '3.147.23'
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("E15").Value = '  -2.05%  '
$ws.Range("D16").Value = '72.325.77'
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("E17").Value = '  -2.27%  '
$ws.Range("D18").Value = '2.659.49'
$ws.Range("E18").Value = '  +0.79%  '
$ws.Range("D19").Value = '12.42'
$ws.Range("E19").Value = '  +5.53%  '
$ws.Range("D20").Value = '370.53'
$ws.Range("E20").Value = '  -3.66%  '
$ws.Range("D21").Value = '7.21'
$ws.Range("E21").Value = '  -9.26%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = '2.09'
$ws.Range("E23").Value = '  +1.28%  '
$ws.Range("D24").Value = '71.93'
$ws.Range("E24").Value = '  -3.14%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").Value = '  -2.23%  '
$ws.Range("D27").Value = '9.86'
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").Value = '0.0₃0974'
$ws.Range("E30").Value = '  +1.32%  '
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("D32").Value = '495.37'
$ws.Range("E32").Value = '  -4.79%  '
$ws.Range("E33").Value = '  -2.74%  '
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").Value = '162.06'
$ws.Range("E36").Value = '  -0.80%  '
$ws.Range("D37").Value = '19.50'
$ws.Range("E37").Value = '  +0.29%  '
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("D39").Value = '18.93'
$ws.Range("E39").Value = '  -0.98%  '
$ws.Range("E40").Value = '  -2.47%  '
$ws.Range("E41").Value = '  -5.34%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").Value = '2.59'
$ws.Range("E43").Value = '  +0.38%  '
$ws.Range("D44").Value = '5.00'
$ws.Range("E44").Value = '  -3.83%  '
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("D46").Value = '155.90'
$ws.Range("E46").Value = '  +3.23%  '
$ws.Range("D47").Value = '39.22'
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("D48").Value = '3.73'
$ws.Range("E48").Value = '  +0.59%  '
$ws.Range("E49").Value = '  +1.73%  '
$ws.Range("E50").Value = '  +1.45%  '
$ws.Range("E51").Value = '  -2.72%  '

# Restore the original (default/unstyled) cell style now that the text values
# are locked in, so no stray style indices are left on the cells.
$editRange.Style = "Normal"
